# Fix same action in one node
# - D10 used to show "B2; BX1" (green "B2", black ";", red " BX1").
#   The "BX1" reference is moved out of D10 into J10 (which already showed "BX2"),
#   so D10 should now only show "B2".
# - J10 used to show plain "BX2". It now shows a rich "BX1;BX2" label: a default/black
#   "BX1", a default/black ";" and a red "BX2" - combining the two actions on one node.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pattern")

# --- D10: remove the red " BX1" portion, leaving just "B2" ---
$d10 = $ws.Range("D10")
$d10.Value = "B2"
$d10.Font.Color = 5296274   # RGB(0,176,80) == FF00B050 green, matches the "B2" run

# --- J10: combine "BX1" and "BX2" into a single rich-text label ---
$j10 = $ws.Range("J10")
$j10.Value = "BX1;BX2"
$j10.Font.Color = 0         # default/automatic black for the whole cell first

# Color just the trailing "BX2" part red, like the rest of the sheet's "BX2" markers
$redStart = "BX1;".Length + 1
$redLen = "BX2".Length
$j10.Characters($redStart, $redLen).Font.Color = 255   # RGB(255,0,0) == FFFF0000 red

# Restore the previously selected range to match the new focus area (F10:F13)
$ws.Range("F10:F13").Select()
